$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 22 (Excel shifts rows 22-28 down to 23-29
# and auto-extends the formulas / dimension / merged cells that cross it).
$ws.Rows("22:22").Insert()

# The blank inserted row doesn't carry over the border styling that the
# rows above/below it have; copy it in from row 21 (format-same-as-above,
# which is what Excel does interactively) for the two columns that need it.
$ws.Range("A21").Copy($ws.Range("A22"))
$ws.Range("F21").Copy($ws.Range("F22"))

# Fill in the new task row.
$ws.Range("A22").Value = "Analyse"
$ws.Range("B22").Value = "Système de gestion de base de données"
$ws.Range("C22").Value = "En cours"
$ws.Range("D22").Value = 10
$ws.Range("E22").Value = 9
$ws.Range("F22").Formula = "=E22/D22"

# Move the selection like the saved workbook (cosmetic, but matches the diff).
$ws.Range("H22").Select()

# Ranges that depended on the old last data row (26) need to grow to 27.
$ws.Range("C1:F27").AutoFilter()
